$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 23:05"

# Row 4: Estados Unidos - refreshed stats
$ws.Range("B4").Value = 1662908
$ws.Range("C4").Value = 17814
$ws.Range("D4").Value = 444505
$ws.Range("E4").Value = 1119837
$ws.Range("G4").Value = 919
$ws.Range("H4").Value = 98566

# Row 40: Rumania - refreshed stats
$ws.Range("E40").Value = 5494
$ws.Range("G40").Value = 10
$ws.Range("H40").Value = 1176

# Row 41: Israel - refreshed stats
$ws.Range("D41").Value = 14090
$ws.Range("E41").Value = 2343

# Rows 140/141: Togo and Cabo Verde swap order (Togo now listed first)
# and both get refreshed stats.
$ws.Range("A140").Value = "Togo"
$ws.Range("B140").Value = 373
$ws.Range("C140").Value = 10
$ws.Range("D140").Value = 133
$ws.Range("E140").Value = 228
$ws.Range("H140").Value = 12

$ws.Range("A141").Value = "Cabo Verde"
$ws.Range("B141").Value = 371
$ws.Range("C141").Value = 9
$ws.Range("D141").Value = 142
$ws.Range("E141").Value = 226
$ws.Range("H141").Value = 3
